# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2310
#   *_new  -> *_FV2404
# Then wrap the data range in an Excel Table ("Table1") and freeze the
# header row, mirroring the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the 21 header cells (row 1, columns A:U) -------------------
$headers = @(
  "Segmentname_FV2310", "Segmentgruppe_FV2310", "Segment_FV2310", "Datenelement_FV2310", "Segment ID_FV2310",
  "Code_FV2310", "Qualifier_FV2310", "Beschreibung_FV2310", "Bedingungsausdruck_FV2310", "Bedingung_FV2310",
  "diff",
  "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404", "Segment ID_FV2404",
  "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404", "Bedingungsausdruck_FV2404", "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the used range A1:U79 into an Excel Table ("Table1") ---------
$tableRange = $ws.Range("A1:U79")
$table = $ws.ListObjects.Add(1, $tableRange, $false, 1, $null)
$table.Name = "Table1"

# --- 3. Freeze the header row (split below row 1) --------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
